$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.176261666666667
$ws.Cells.Item(2, 8).Value = 3.528785
$ws.Cells.Item(2, 9).Value = 0.003383077204500727
$ws.Cells.Item(2, 10).Value = 0.003383077204500727
$ws.Cells.Item(2, 13).Value = 24.576554
$ws.Cells.Item(2, 14).Value = 73.729662
$ws.Cells.Item(2, 15).Value = 0.07553767049546639
$ws.Cells.Item(2, 16).Value = 0.07553767049546638
$ws.Cells.Item(2, 17).Value = 28.90845836896334
$ws.Cells.Item(2, 18).Value = 260.17612532067
$ws.Cells.Item(2, 19).Value = 0.0002555497711342995
$ws.Cells.Item(2, 20).Value = 0.0002555497711342994
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.176261666666667
$ws.Cells.Item(3, 8).Value = 3.528785
$ws.Cells.Item(3, 9).Value = 0.003383077204500727
$ws.Cells.Item(3, 10).Value = 0.003383077204500727
$ws.Cells.Item(3, 15).Value = 0.359764849016532
$ws.Cells.Item(3, 16).Value = 0.359764849016532
$ws.Cells.Item(3, 17).Value = 137.6829215435628
$ws.Cells.Item(3, 18).Value = 1239.146293892065
$ws.Cells.Item(3, 19).Value = 0.001217112259688475
$ws.Cells.Item(3, 20).Value = 0.001217112259688475
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.176261666666667
$ws.Cells.Item(4, 8).Value = 3.528785
$ws.Cells.Item(4, 9).Value = 0.003383077204500727
$ws.Cells.Item(4, 10).Value = 0.003383077204500727
$ws.Cells.Item(4, 13).Value = 55.68784966666667
$ws.Cells.Item(4, 14).Value = 167.063549
$ws.Cells.Item(4, 15).Value = 0.1711603033819035
$ws.Cells.Item(4, 16).Value = 0.1711603033819035
$ws.Cells.Item(4, 17).Value = 65.50348286199612
$ws.Cells.Item(4, 18).Value = 589.5313457579651
$ws.Cells.Item(4, 19).Value = 0.0005790485206867464
$ws.Cells.Item(4, 20).Value = 0.0005790485206867464
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.176261666666667
$ws.Cells.Item(5, 8).Value = 3.528785
$ws.Cells.Item(5, 9).Value = 0.003383077204500727
$ws.Cells.Item(5, 10).Value = 0.003383077204500727
$ws.Cells.Item(5, 13).Value = 128.0392633333333
$ws.Cells.Item(5, 14).Value = 384.11779
$ws.Cells.Item(5, 15).Value = 0.3935371771060981
$ws.Cells.Item(5, 16).Value = 0.3935371771060981
$ws.Cells.Item(5, 17).Value = 150.6076772872389
$ws.Cells.Item(5, 18).Value = 1355.46909558515
$ws.Cells.Item(5, 19).Value = 0.001331366652991206
$ws.Cells.Item(5, 20).Value = 0.001331366652991206
$ws.Cells.Item(6, 9).Value = 0.9815272193474828
$ws.Cells.Item(6, 10).Value = 0.9815272193474829
$ws.Cells.Item(6, 13).Value = 24.576554
$ws.Cells.Item(6, 14).Value = 73.729662
$ws.Cells.Item(6, 15).Value = 0.07553767049546639
$ws.Cells.Item(6, 16).Value = 0.07553767049546638
$ws.Cells.Item(6, 17).Value = 8387.168557892412
$ws.Cells.Item(6, 18).Value = 75484.5170210317
$ws.Cells.Item(6, 19).Value = 0.07414227967740152
$ws.Cells.Item(6, 20).Value = 0.07414227967740152
$ws.Cells.Item(7, 9).Value = 0.9815272193474828
$ws.Cells.Item(7, 10).Value = 0.9815272193474829
$ws.Cells.Item(7, 15).Value = 0.359764849016532
$ws.Cells.Item(7, 16).Value = 0.359764849016532
$ws.Cells.Item(7, 17).Value = 39945.74376088904
$ws.Cells.Item(7, 18).Value = 359511.6938480013
$ws.Cells.Item(7, 19).Value = 0.3531189918741636
$ws.Cells.Item(7, 20).Value = 0.3531189918741637
$ws.Cells.Item(8, 9).Value = 0.9815272193474828
$ws.Cells.Item(8, 10).Value = 0.9815272193474829
$ws.Cells.Item(8, 13).Value = 55.68784966666667
$ws.Cells.Item(8, 14).Value = 167.063549
$ws.Cells.Item(8, 15).Value = 0.1711603033819035
$ws.Cells.Item(8, 16).Value = 0.1711603033819035
$ws.Cells.Item(8, 17).Value = 19004.42925321858
$ws.Cells.Item(8, 18).Value = 171039.8632789672
$ws.Cells.Item(8, 19).Value = 0.1679984966411113
$ws.Cells.Item(8, 20).Value = 0.1679984966411113
$ws.Cells.Item(9, 9).Value = 0.9815272193474828
$ws.Cells.Item(9, 10).Value = 0.9815272193474829
$ws.Cells.Item(9, 13).Value = 128.0392633333333
$ws.Cells.Item(9, 14).Value = 384.11779
$ws.Cells.Item(9, 15).Value = 0.3935371771060981
$ws.Cells.Item(9, 16).Value = 0.3935371771060981
$ws.Cells.Item(9, 17).Value = 43695.5841573656
$ws.Cells.Item(9, 18).Value = 393260.2574162903
$ws.Cells.Item(9, 19).Value = 0.3862674511548063
$ws.Cells.Item(9, 20).Value = 0.3862674511548064
$ws.Cells.Item(10, 7).Value = 2.611920666666667
$ws.Cells.Item(10, 8).Value = 7.835762
$ws.Cells.Item(10, 9).Value = 0.007512213921248538
$ws.Cells.Item(10, 10).Value = 0.007512213921248538
$ws.Cells.Item(10, 13).Value = 24.576554
$ws.Cells.Item(10, 14).Value = 73.729662
$ws.Cells.Item(10, 15).Value = 0.07553767049546639
$ws.Cells.Item(10, 16).Value = 0.07553767049546638
$ws.Cells.Item(10, 17).Value = 64.19200930804934
$ws.Cells.Item(10, 18).Value = 577.7280837724441
$ws.Cells.Item(10, 19).Value = 0.0005674551398747276
$ws.Cells.Item(10, 20).Value = 0.0005674551398747275
$ws.Cells.Item(11, 7).Value = 2.611920666666667
$ws.Cells.Item(11, 8).Value = 7.835762
$ws.Cells.Item(11, 9).Value = 0.007512213921248538
$ws.Cells.Item(11, 10).Value = 0.007512213921248538
$ws.Cells.Item(11, 15).Value = 0.359764849016532
$ws.Cells.Item(11, 16).Value = 0.359764849016532
$ws.Cells.Item(11, 17).Value = 305.7286303019398
$ws.Cells.Item(11, 18).Value = 2751.557672717458
$ws.Cells.Item(11, 19).Value = 0.00270263050715787
$ws.Cells.Item(11, 20).Value = 0.00270263050715787
$ws.Cells.Item(12, 7).Value = 2.611920666666667
$ws.Cells.Item(12, 8).Value = 7.835762
$ws.Cells.Item(12, 9).Value = 0.007512213921248538
$ws.Cells.Item(12, 10).Value = 0.007512213921248538
$ws.Cells.Item(12, 13).Value = 55.68784966666667
$ws.Cells.Item(12, 14).Value = 167.063549
$ws.Cells.Item(12, 15).Value = 0.1711603033819035
$ws.Cells.Item(12, 16).Value = 0.1711603033819035
$ws.Cells.Item(12, 17).Value = 145.4522454265931
$ws.Cells.Item(12, 18).Value = 1309.070208839338
$ws.Cells.Item(12, 19).Value = 0.001285792813830659
$ws.Cells.Item(12, 20).Value = 0.001285792813830659
$ws.Cells.Item(13, 7).Value = 2.611920666666667
$ws.Cells.Item(13, 8).Value = 7.835762
$ws.Cells.Item(13, 9).Value = 0.007512213921248538
$ws.Cells.Item(13, 10).Value = 0.007512213921248538
$ws.Cells.Item(13, 13).Value = 128.0392633333333
$ws.Cells.Item(13, 14).Value = 384.11779
$ws.Cells.Item(13, 15).Value = 0.3935371771060981
$ws.Cells.Item(13, 16).Value = 0.3935371771060981
$ws.Cells.Item(13, 17).Value = 334.4283980451089
$ws.Cells.Item(13, 18).Value = 3009.85558240598
$ws.Cells.Item(13, 19).Value = 0.002956335460385282
$ws.Cells.Item(13, 20).Value = 0.002956335460385282
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 2.634616333333333
$ws.Cells.Item(14, 8).Value = 7.903849
$ws.Cells.Item(14, 9).Value = 0.007577489526767955
$ws.Cells.Item(14, 10).Value = 0.007577489526767956
$ws.Cells.Item(14, 13).Value = 24.576554
$ws.Cells.Item(14, 14).Value = 73.729662
$ws.Cells.Item(14, 15).Value = 0.07553767049546639
$ws.Cells.Item(14, 16).Value = 0.07553767049546638
$ws.Cells.Item(14, 17).Value = 64.74979058544866
$ws.Cells.Item(14, 18).Value = 582.7481152690381
$ws.Cells.Item(14, 19).Value = 0.0005723859070558454
$ws.Cells.Item(14, 20).Value = 0.0005723859070558454
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 2.634616333333333
$ws.Cells.Item(15, 8).Value = 7.903849
$ws.Cells.Item(15, 9).Value = 0.007577489526767955
$ws.Cells.Item(15, 10).Value = 0.007577489526767956
$ws.Cells.Item(15, 15).Value = 0.359764849016532
$ws.Cells.Item(15, 16).Value = 0.359764849016532
$ws.Cells.Item(15, 17).Value = 308.3851869012046
$ws.Cells.Item(15, 18).Value = 2775.466682110841
$ws.Cells.Item(15, 19).Value = 0.002726114375522026
$ws.Cells.Item(15, 20).Value = 0.002726114375522026
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 2.634616333333333
$ws.Cells.Item(16, 8).Value = 7.903849
$ws.Cells.Item(16, 9).Value = 0.007577489526767955
$ws.Cells.Item(16, 10).Value = 0.007577489526767956
$ws.Cells.Item(16, 13).Value = 55.68784966666667
$ws.Cells.Item(16, 14).Value = 167.063549
$ws.Cells.Item(16, 15).Value = 0.1711603033819035
$ws.Cells.Item(16, 16).Value = 0.1711603033819035
$ws.Cells.Item(16, 17).Value = 146.7161183000112
$ws.Cells.Item(16, 18).Value = 1320.445064700101
$ws.Cells.Item(16, 19).Value = 0.0012969654062748
$ws.Cells.Item(16, 20).Value = 0.0012969654062748
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 2.634616333333333
$ws.Cells.Item(17, 8).Value = 7.903849
$ws.Cells.Item(17, 9).Value = 0.007577489526767955
$ws.Cells.Item(17, 10).Value = 0.007577489526767956
$ws.Cells.Item(17, 13).Value = 128.0392633333333
$ws.Cells.Item(17, 14).Value = 384.11779
$ws.Cells.Item(17, 15).Value = 0.3935371771060981
$ws.Cells.Item(17, 16).Value = 0.3935371771060981
$ws.Cells.Item(17, 17).Value = 337.3343344859678
$ws.Cells.Item(17, 18).Value = 3036.00901037371
$ws.Cells.Item(17, 19).Value = 0.002982023837915284
$ws.Cells.Item(17, 20).Value = 0.002982023837915285

Write-Output "Applied 192 cell updates"